# The deck's single slide master currently carries the "Integral" theme
# colour scheme. The target state recolours it to the stock "Office Theme"
# palette instead (the colours that, before this edit, were only used by
# the - otherwise practically invisible - notes master).
#
# PowerPoint's theme colours are edited through the Theme's
# ThemeColorScheme collection (Design > Variants > Colors > Customize
# Colors...), 12 slots in clrScheme document order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. Each item's .RGB takes a standard VBA
# RGB() long (0x00BBGGRR).

function RGB([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1      000000
$cs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$cs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      44546A
$cs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$cs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$cs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  ED7D31
$cs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$cs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  FFC000
$cs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  4472C4
$cs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6  70AD47
$cs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink    0563C1
$cs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink 954F72
